# Apply updated cryptocurrency data to sheet1 (cell values only; no structural changes).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.453.60"
$ws.Range("E2").Value = "  -0.85%  "
$ws.Range("D3").Value = "1.572.79"
$ws.Range("E3").Value = "  -0.75%  "
$ws.Range("E4").Value = "  -0.35%  "
$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "208.13"
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = "  +0.58%  "
$origStyle = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.500"
$ws.Range("D6").Style = $origStyle
$ws.Range("E6").Value = "  -0.77%  "
$ws.Range("E7").Value = "  -0.33%  "
$origStyle = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.28"
$ws.Range("D8").Style = $origStyle
$ws.Range("E8").Value = "  +0.55%  "
$origStyle = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.249"
$ws.Range("D9").Style = $origStyle
$ws.Range("E9").Value = "  -0.85%  "
$origStyle = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0596"
$ws.Range("D10").Style = $origStyle
$ws.Range("E10").Value = "  +1.24%  "
$ws.Range("E11").Value = "  -0.27%  "
$ws.Range("D12").Value = "1.795.78"
$ws.Range("E12").Value = "  -0.86%  "
$ws.Range("D13").Value = "1.568.03"
$ws.Range("E13").Value = "  -1.03%  "
$ws.Range("E14").Value = "  -1.26%  "
$origStyle = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.520"
$ws.Range("D15").Style = $origStyle
$ws.Range("E15").Value = "  -1.78%  "
$origStyle = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.61"
$ws.Range("D16").Style = $origStyle
$ws.Range("E16").Value = "  +0.54%  "
$ws.Range("D17").Value = "27.439.73"
$ws.Range("E17").Value = "  -0.76%  "
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").Value = "0.0₃0694"
$ws.Range("E18").Value = "  +0.12%  "
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$origStyle = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.32"
$ws.Range("D19").Style = $origStyle
$ws.Range("E19").Value = "  -0.05%  "
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$origStyle = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "213.00"
$ws.Range("D20").Style = $origStyle
$ws.Range("E20").Value = "  -3.50%  "
$ws.Range("E21").Value = "  -0.35%  "
$origStyle = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.13"
$ws.Range("D22").Style = $origStyle
$ws.Range("E22").Value = "  -0.13%  "
$origStyle = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.55"
$ws.Range("D23").Style = $origStyle
$ws.Range("E23").Value = "  -0.38%  "
$origStyle = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.00"
$ws.Range("D24").Style = $origStyle
$ws.Range("E24").Value = "  +0.85%  "
$origStyle = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.23"
$ws.Range("D25").Style = $origStyle
$ws.Range("E25").Value = "  -0.53%  "
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$origStyle = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.69"
$ws.Range("D26").Style = $origStyle
$ws.Range("E26").Value = "  -2.75%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$origStyle = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "14.97"
$ws.Range("D27").Style = $origStyle
$ws.Range("E27").Value = "  -0.77%  "
$ws.Range("B28").Value = "Stellar"
$ws.Range("C28").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$origStyle = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.105"
$ws.Range("D28").Style = $origStyle
$ws.Range("E28").Value = "  -1.05%  "
$ws.Range("B29").Value = "BinanceUSD"
$ws.Range("C29").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$origStyle = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = $origStyle
$ws.Range("E29").Value = "  -0.31%  "
$ws.Range("E30").Value = "  -0.06%  "
$origStyle = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0472"
$ws.Range("D31").Style = $origStyle
$ws.Range("E31").Value = "  +0.86%  "
$ws.Range("E32").Value = "  -0.87%  "
$ws.Range("D33").Value = "1.390.46"
$ws.Range("E33").Value = "  +1.64%  "
$ws.Range("E34").Value = "  +1.35%  "
$origStyle = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.55"
$ws.Range("D35").Style = $origStyle
$ws.Range("E35").Value = "  +1.54%  "
$origStyle = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.950"
$ws.Range("D36").Style = $origStyle
$ws.Range("E36").Value = "  -2.94%  "
$origStyle = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.29"
$ws.Range("D37").Style = $origStyle
$ws.Range("E37").Value = "  -1.17%  "
$ws.Range("E38").Value = "  -0.12%  "
$origStyle = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.534"
$ws.Range("D39").Style = $origStyle
$ws.Range("E39").Value = "  -0.44%  "
$ws.Range("E40").Value = "  +0.55%  "
$ws.Range("E41").Value = "  -0.31%  "
$origStyle = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.992"
$ws.Range("D42").Style = $origStyle
$ws.Range("E42").Value = "  +1.08%  "
$ws.Range("E43").Value = "  +3.23%  "
$origStyle = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "64.28"
$ws.Range("D44").Style = $origStyle
$ws.Range("E44").Value = "  +1.38%  "
$ws.Range("E45").Value = "  +0.01%  "
$origStyle = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.27"
$ws.Range("D46").Style = $origStyle
$ws.Range("E46").Value = "  +0.94%  "
$ws.Range("D47").Value = "1.707.44"
$origStyle = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "85.71"
$ws.Range("D48").Style = $origStyle
$ws.Range("E48").Value = "  -2.97%  "
$ws.Range("D49").Value = "0.0₇0989"
$ws.Range("E49").Value = "  -1.73%  "
$origStyle = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0957"
$ws.Range("D50").Style = $origStyle
$ws.Range("E50").Value = "  -1.15%  "
$origStyle = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0494"
$ws.Range("D51").Style = $origStyle
$ws.Range("E51").Value = "  -0.84%  "
